# Auto-generated edit script for LOT2059.xlsx
# Applies the shared-strings/content realignment + two new trailing rows
# described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Clear cells that must no longer hold any value (incl. removing their cell record) ---
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B19").Clear()
$ws.Range("C19").Clear()

# --- 2) Write final text for every populated cell A1:C25 ---
$ws.Range("B1").Value2 = "Ementa atual:"
$ws.Range("C1").Value2 = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value2 = "LOT2059"
$ws.Range("C2").Value2 = "LOT2059"
$ws.Range("A3").Value2 = "Nome:"
$ws.Range("B3").Value2 = " Química Orgânica Fundamental"
$ws.Range("C3").Value2 = " Química Orgânica Fundamental"
$ws.Range("A4").Value2 = "Name:"
$ws.Range("B4").Value2 = "Basics of Organic Chemistry"
$ws.Range("C4").Value2 = "Basics of Organic Chemistry"
$ws.Range("A5").Value2 = "Créditos-aula:"
$ws.Range("B5").Value2 = "4"
$ws.Range("C5").Value2 = "4"
$ws.Range("A6").Value2 = "Créditos-trabalho"
$ws.Range("B6").Value2 = "0"
$ws.Range("C6").Value2 = "0"
$ws.Range("A7").Value2 = "Carga horária:"
$ws.Range("B7").Value2 = "60 h"
$ws.Range("C7").Value2 = "60 h"
$ws.Range("A8").Value2 = "Ativação:"
$ws.Range("B8").Value2 = "01/01/2017"
$ws.Range("C8").Value2 = "01/01/2017"
$ws.Range("A9").Value2 = "Semestre ideal:"
$ws.Range("B9").Value2 = "EB-3"
$ws.Range("C9").Value2 = "EB-3"
$ws.Range("A10").Value2 = "Objetivos:"
$ws.Range("B10").Value2 = "Introdução teórica da Química Orgânica aos estudantes de Engenharia Bioquímica abordando de forma sistematizada: a) a relação entre a estrutura das moléculas orgânicas e suas propriedades físico-químicas; b) a reatividade das moléculas orgânicas em função do tipo de grupo funcional que carregam e; c) as principais vias de reações entre moléculas orgânicas."
$ws.Range("C10").Value2 = "Introdução teórica da Química Orgânica aos estudantes de Engenharia Bioquímica abordando de forma sistematizada: a) a relação entre a estrutura das moléculas orgânicas e suas propriedades físico-químicas; b) a reatividade das moléculas orgânicas em função do tipo de grupo funcional que carregam e; c) as principais vias de reações entre moléculas orgânicas."
$ws.Range("A11").Value2 = "Objectives:"
$ws.Range("A12").Value2 = "Docentes responsáveis:"
$ws.Range("B13").Value2 = "2143261 - André Luis Ferraz"
$ws.Range("C13").Value2 = "2143261 - André Luis Ferraz"
$ws.Range("B14").Value2 = "3380737 - Flávio Teixeira da Silva"
$ws.Range("C14").Value2 = "3380737 - Flávio Teixeira da Silva"
$ws.Range("A15").Value2 = "Programa resumido:"
$ws.Range("B15").Value2 = "A disciplina abordará os fundamentos da química orgânica que darão apoio às disciplinas subsequentes na área de bioquímica, biologia molecular, polímeros e química de biomassa. A abordagem teórica dará subsídios ao aluno para interpretar as propriedades e a reatividade das moléculas orgânicas desde um ponto de vista estrutural."
$ws.Range("C15").Value2 = "A disciplina abordará os fundamentos da química orgânica que darão apoio às disciplinas subsequentes na área de bioquímica, biologia molecular, polímeros e química de biomassa. A abordagem teórica dará subsídios ao aluno para interpretar as propriedades e a reatividade das moléculas orgânicas desde um ponto de vista estrutural."
$ws.Range("A16").Value2 = "Short syllabus:"
$ws.Range("A17").Value2 = "Programa:"
$ws.Range("B17").Value2 = "Estrutura versus propriedades físico-química dos: hidrocarbonetos, compostos com grupos funcionais formados por ligações simples, o grupo carbonila e seus compostos derivados, grupos funcionais que contém hetero-átomos, benzeno e aromaticidade;Esteroquímica; Reações químicas de compostos orgânicos: reações de alcenos e alcinos (adições à dupla ligação); reações de compostos aromáticos (substituição nucleofílica em aromáticos); reações de compostos orgânicos halogenados (substituição nucleofílica e eliminação); reações de álcoois, fenóis e éteres; reações de aldeídos e cetonas (adições em compostos carbonílicos); reações de ácidos carboxílicos e derivados."
$ws.Range("C17").Value2 = "Estrutura versus propriedades físico-química dos: hidrocarbonetos, compostos com grupos funcionais formados por ligações simples, o grupo carbonila e seus compostos derivados, grupos funcionais que contém hetero-átomos, benzeno e aromaticidade;Esteroquímica; Reações químicas de compostos orgânicos: reações de alcenos e alcinos (adições à dupla ligação); reações de compostos aromáticos (substituição nucleofílica em aromáticos); reações de compostos orgânicos halogenados (substituição nucleofílica e eliminação); reações de álcoois, fenóis e éteres; reações de aldeídos e cetonas (adições em compostos carbonílicos); reações de ácidos carboxílicos e derivados."
$ws.Range("A18").Value2 = "Syllabus:"
$ws.Range("A19").Value2 = "Avaliação:"
$ws.Range("A20").Value2 = "Método:"
$ws.Range("B20").Value2 = "A avaliação será feita por meio de provas escritas (P1 e P2). Critério"
$ws.Range("C20").Value2 = "A avaliação será feita por meio de provas escritas (P1 e P2). Critério"
$ws.Range("A21").Value2 = "Critério:"
$ws.Range("B21").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = [(P1x1) + (P2x2)]/3"
$ws.Range("C21").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = [(P1x1) + (P2x2)]/3"
$ws.Range("A22").Value2 = "Norma de recuperação:"
$ws.Range("B22").Value2 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C22").Value2 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("A23").Value2 = "Bibliografia:"
$ws.Range("B23").Value2 = "SOLOMONS, T.W.G., FRYHLE, C.B. Química Orgânica 1 e 2. 10ª Edição, Rio de Janeiro, LTC Editora, 2012.BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo, Pearson Prentice Hall, 2006.ALLINGER, N.L. Química Orgânica, 2ª Edição, Rio de Janeiro, Guanabara Dois, 1976."
$ws.Range("C23").Value2 = "SOLOMONS, T.W.G., FRYHLE, C.B. Química Orgânica 1 e 2. 10ª Edição, Rio de Janeiro, LTC Editora, 2012.BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo, Pearson Prentice Hall, 2006.ALLINGER, N.L. Química Orgânica, 2ª Edição, Rio de Janeiro, Guanabara Dois, 1976."
$ws.Range("A24").Value2 = "Requisitos:"
$ws.Range("B25").Value2 = "LOQ4073 -  Química Geral II  (Requisito fraco)`n"
$ws.Range("C25").Value2 = "LOQ4073 -  Química Geral II  (Requisito fraco)`n"

# --- 3) Newly-introduced cells need their formatting copied from a same-column sibling ---
#        (column A -> style of A9, column B -> style of B9, column C -> style of C9)
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 4) Row heights that differ from the sheet default ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30

Write-Host "Edit complete"
